# Actualizacion de tarea completada
# Marks the "Validar fechas..." task (row 27) as completed and adds a new
# pending task row for CUIT validation message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 27: task "Validar fechas - permite cargar fechas q no existen" is now
# finished -> change status cell from text "en proceso" to a 100% (done) value.
$estado = $ws.Cells.Item(27, 3)
$estado.Value = 1
$estado.NumberFormat = "0%"

# New task row (29): add the new pending task description.
$ws.Cells.Item(29, 1).Value = "Validacion de cuit para mostrar mensaje correcto"

# New empty row (30) with the underlined "separator" style used elsewhere
# in the sheet (e.g. C13, D11).
$ws.Cells.Item(30, 1).Font.Underline = 2

# Reset the view: scroll back to the top-left and select cell B4.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B4").Select()
